$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.36818657311801
$ws.Range("D2").Value = 0.1540160869351439
$ws.Range("E2").Value = 0.8520184962022483
$ws.Range("F2").Value = 3.155443518393298
$ws.Range("G2").Value = 0.002431486114854218
$ws.Range("I2").Value = 0.8585112567489368
$ws.Range("L2").Value = 0.8819045890216728
$ws.Range("N2").Value = 1.314705007570993

$ws.Range("B3").Value = 1.278937653627281
$ws.Range("D3").Value = 0.1434341466781888
$ws.Range("E3").Value = 0.7408576144380277
$ws.Range("F3").Value = 2.937355339469775
$ws.Range("G3").Value = 0.00244324169056885
$ws.Range("I3").Value = 0.8725010023283737
$ws.Range("L3").Value = 0.7992595441351682
$ws.Range("N3").Value = 1.331890729450329

$ws.Range("B4").Value = 1.22480706340275
$ws.Range("D4").Value = 0.1371046108474445
$ws.Range("E4").Value = 0.6727087933125517
$ws.Range("F4").Value = 2.806495658123339
$ws.Range("G4").Value = 0.002450810248189074
$ws.Range("I4").Value = 0.8817711296233277
$ws.Range("L4").Value = 0.7490527470005475
$ws.Range("N4").Value = 1.343066689298631

$ws.Range("B5").Value = 1.202914322426921
$ws.Range("D5").Value = 0.1345660214885243
$ws.Range("E5").Value = 0.6449553458229929
$ws.Range("F5").Value = 2.753910692077255
$ws.Range("G5").Value = 0.002453983147160446
$ws.Range("I5").Value = 0.8857187580388164
$ws.Range("L5").Value = 0.7287232051021135
$ws.Range("N5").Value = 1.347778055795445

$ws.Range("B6").Value = 1.199289007401092
$ws.Range("D6").Value = 0.1341469099328236
$ws.Range("E6").Value = 0.6403477348806774
$ws.Range("F6").Value = 2.745223062781008
$ws.Range("G6").Value = 0.002454515373107096
$ws.Range("I6").Value = 0.886384492879781
$ws.Range("L6").Value = 0.7253552150854148
$ws.Range("N6").Value = 1.348569867440474

$ws.Range("B7").Value = 1.224511141318544
$ws.Range("D7").Value = 0.1370702114505917
$ws.Range("E7").Value = 0.6723344398210287
$ws.Range("F7").Value = 2.805783509330951
$ws.Range("G7").Value = 0.002450852679441275
$ws.Range("I7").Value = 0.881823681936396
$ws.Range("L7").Value = 0.7487780553476
$ws.Range("N7").Value = 1.343129592178023

$ws.Range("B8").Value = 1.337273388075573
$ws.Range("D8").Value = 0.150331716371042
$ws.Range("E8").Value = 0.8136633360409462
$ws.Range("F8").Value = 3.079599026445351
$ws.Range("G8").Value = 0.002435466971560117
$ws.Range("I8").Value = 0.8631930394434377
$ws.Range("L8").Value = 0.853293896723585
$ws.Range("N8").Value = 1.32050127969525

$ws.Range("B9").Value = 1.563821209404978
$ws.Range("D9").Value = 0.1777372727593161
$ws.Range("E9").Value = 1.09204304113976
$ws.Range("F9").Value = 3.641904763410452
$ws.Range("G9").Value = 0.002408054711367846
$ws.Range("I9").Value = 0.8321071925946839
$ws.Range("L9").Value = 1.06275854833649
$ws.Range("N9").Value = 1.281068576760745

$ws.Range("B10").Value = 1.733756263940222
$ws.Range("D10").Value = 0.1988241522265071
$ws.Range("E10").Value = 1.297916314104526
$ws.Range("F10").Value = 4.07219582920635
$ws.Range("G10").Value = 0.002389564823115581
$ws.Range("I10").Value = 0.8126606501025435
$ws.Range("L10").Value = 1.219763103422963
$ws.Range("N10").Value = 1.255098126836629

$ws.Range("B11").Value = 1.8118653800326
$ws.Range("D11").Value = 0.2086461498382732
$ws.Range("E11").Value = 1.392001735383019
$ws.Range("F11").Value = 4.27206068026635
$ws.Range("G11").Value = 0.002381504524631012
$ws.Range("I11").Value = 0.8045662546886376
$ws.Range("L11").Value = 1.291946544248276
$ws.Range("N11").Value = 1.243932866996765

$ws.Range("B12").Value = 1.841562334198102
$ws.Range("D12").Value = 0.2124003429671291
$ws.Range("E12").Value = 1.427702882526347
$ws.Range("F12").Value = 4.348369373916398
$ws.Range("G12").Value = 0.00237850218745761
$ws.Range("I12").Value = 0.8016106063863191
$ws.Range("L12").Value = 1.319396719696329
$ws.Range("N12").Value = 1.239798042791932

$ws.Range("B13").Value = 1.835161234072018
$ws.Range("D13").Value = 0.2115902322521208
$ws.Range("E13").Value = 1.42001057635224
$ws.Range("F13").Value = 4.331906659552487
$ws.Range("G13").Value = 0.002379146582510851
$ws.Range("I13").Value = 0.8022422644013432
$ws.Range("L13").Value = 1.313479575310112
$ws.Range("N13").Value = 1.24068440647644

$ws.Range("B14").Value = 1.814306169177485
$ws.Range("D14").Value = 0.2089543011674948
$ws.Range("E14").Value = 1.394937362878323
$ws.Range("F14").Value = 4.278325949180442
$ws.Range("G14").Value = 0.002381256523086456
$ws.Range("I14").Value = 0.8043208885813158
$ws.Range("L14").Value = 1.294202526920515
$ws.Range("N14").Value = 1.243590824739712

$ws.Range("B15").Value = 1.80154739294494
$ws.Range("D15").Value = 0.2073443063311231
$ws.Range("E15").Value = 1.379589136129908
$ws.Range("F15").Value = 4.245588454954827
$ws.Range("G15").Value = 0.002382555408021645
$ws.Range("I15").Value = 0.8056084111939228
$ws.Range("L15").Value = 1.282410062291035
$ws.Range("N15").Value = 1.245383227260156

$ws.Range("B16").Value = 1.728668097164814
$ws.Range("D16").Value = 0.1981870331531184
$ws.Range("E16").Value = 1.291777232434811
$ws.Range("F16").Value = 4.059219812346441
$ws.Range("G16").Value = 0.002390098597604575
$ws.Range("I16").Value = 0.8132048872225397
$ws.Range("L16").Value = 1.215061606612096
$ws.Range("N16").Value = 1.255840851462381

$ws.Range("B17").Value = 1.684167227023863
$ws.Range("D17").Value = 0.1926293900930602
$ws.Range("E17").Value = 1.238025856290051
$ws.Range("F17").Value = 3.945966683786878
$ws.Range("G17").Value = 0.002394815586503698
$ws.Range("I17").Value = 0.818058649280637
$ws.Range("L17").Value = 1.173944973992093
$ws.Range("N17").Value = 1.262422381872007

$ws.Range("B18").Value = 1.658647067868742
$ws.Range("D18").Value = 0.1894542672562523
$ws.Range("E18").Value = 1.207149535909508
$ws.Range("F18").Value = 3.881212647439952
$ws.Range("G18").Value = 0.002397561728925534
$ws.Range("I18").Value = 0.8209210947100019
$ws.Range("L18").Value = 1.150367032868814
$ws.Range("N18").Value = 1.266268978964987

$ws.Range("B19").Value = 1.650019269766631
$ws.Range("D19").Value = 0.1883828611274225
$ws.Range("E19").Value = 1.196701899719585
$ws.Range("F19").Value = 3.859353432162663
$ws.Range("G19").Value = 0.00239849721882697
$ws.Range("I19").Value = 0.8219023686655191
$ws.Range("L19").Value = 1.142395994598303
$ws.Range("N19").Value = 1.267581863081134

$ws.Range("B20").Value = 1.688896578749564
$ws.Range("D20").Value = 0.1932187739415383
$ws.Range("E20").Value = 1.243743559955732
$ws.Range("F20").Value = 3.95798248254664
$ws.Range("G20").Value = 0.0023943100376036
$ws.Range("I20").Value = 0.8175346318977077
$ws.Range("L20").Value = 1.178314489121362
$ws.Range("N20").Value = 1.261715446356519

$ws.Range("B21").Value = 1.820428560439609
$ws.Range("D21").Value = 0.2097275783388852
$ws.Range("E21").Value = 1.402299902825945
$ws.Range("F21").Value = 4.294046706943618
$ws.Range("G21").Value = 0.002380635431732195
$ws.Range("I21").Value = 0.8037073625128386
$ws.Range("L21").Value = 1.299861466359744
$ws.Range("N21").Value = 1.242734609939532

$ws.Range("B22").Value = 1.907085842569586
$ws.Range("D22").Value = 0.2207208578450093
$ws.Range("E22").Value = 1.506356699334191
$ws.Range("F22").Value = 4.517337797447567
$ws.Range("G22").Value = 0.002371989024857895
$ws.Range("I22").Value = 0.7953095053600023
$ws.Range("L22").Value = 1.379977981168793
$ws.Range("N22").Value = 1.230872849642886

$ws.Range("B23").Value = 1.860770733859852
$ws.Range("D23").Value = 0.2148342710607949
$ws.Range("E23").Value = 1.450776499310024
$ws.Range("F23").Value = 4.397818203655788
$ws.Range("G23").Value = 0.002376577348701716
$ws.Range("I23").Value = 0.7997326541214349
$ws.Range("L23").Value = 1.337154025666621
$ws.Range("N23").Value = 1.237154004991105

$ws.Range("B24").Value = 1.68675823976406
$ws.Range("D24").Value = 0.1929522512194239
$ws.Range("E24").Value = 1.241158505765895
$ws.Range("F24").Value = 3.952549031253398
$ws.Range("G24").Value = 0.002394538489274391
$ws.Range("I24").Value = 0.8177713161443378
$ws.Range("L24").Value = 1.176338843186159
$ws.Range("N24").Value = 1.262034856327176

$ws.Range("B25").Value = 1.501934918412246
$ws.Range("D25").Value = 0.1701634738331279
$ws.Range("E25").Value = 1.016541552779131
$ws.Range("F25").Value = 3.486898513492719
$ws.Range("G25").Value = 0.002415178372136457
$ws.Range("I25").Value = 0.8399262373305731
$ws.Range("L25").Value = 1.005573370878125
$ws.Range("N25").Value = 1.291208416706894
